$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.26437026758141258
$ws.Range("B1").Value = 0.2637540865426331
$ws.Range("A2").Value = -0.16787880015843015
$ws.Range("B2").Value = 0.16639073009917738
$ws.Range("A3").Value = -0.063445774756303308
$ws.Range("B3").Value = 0.063160602815663935
$ws.Range("A4").Value = -0.12515725109119558
$ws.Range("B4").Value = 0.12460964881125136
$ws.Range("A5").Value = -0.1186096490737425
$ws.Range("B5").Value = 0.11753288574789167
$ws.Range("A6").Value = -0.08636520852316254
$ws.Range("B6").Value = 0.086246930049843407
$ws.Range("A7").Value = -0.066246930374690649
$ws.Range("B7").Value = 0.065980010127395516
$ws.Range("A8").Value = -0.045980010455368259
$ws.Range("B8").Value = 0.045788881434422279
$ws.Range("A9").Value = -0.039788881711423585
$ws.Range("B9").Value = 0.039638147473302965
$ws.Range("A10").Value = -0.033638147753215719
$ws.Range("B10").Value = 0.03361946432315932
$ws.Range("A11").Value = -0.029119464597439304
$ws.Range("B11").Value = 0.029090586758655945
$ws.Range("A12").Value = -0.023090587039431121
$ws.Range("B12").Value = 0.023014367448852102
$ws.Range("A13").Value = -0.017014367732062219
$ws.Range("B13").Value = 0.016998838212175293
$ws.Range("A14").Value = -0.02708635341274146
$ws.Range("B14").Value = 0.027053703203619328
$ws.Range("A15").Value = -0.021053703488550291
$ws.Range("B15").Value = 0.021027995178125813
$ws.Range("A16").Value = -0.015027995464103938
$ws.Range("B16").Value = 0.015004597356626626
$ws.Range("A17").Value = -0.0090045976439609987
$ws.Range("B17").Value = 0.0089999997006815491
$ws.Range("A18").Value = -0.085875286652321137
$ws.Range("B18").Value = 0.08577799093366778
$ws.Range("A19").Value = -0.027097192854809826
$ws.Range("B19").Value = 0.027013941966810506
$ws.Range("A20").Value = -0.018013942234404112
$ws.Range("B20").Value = 0.018004312687219581
$ws.Range("A21").Value = -0.0090043129552084267
$ws.Range("B21").Value = 0.0089999997317127267
$ws.Range("A22").Value = -0.093944955612483838
$ws.Range("B22").Value = 0.093632623691572547
$ws.Range("A23").Value = -0.084632623963790898
$ws.Range("B23").Value = 0.084126502754730659
$ws.Range("A24").Value = -0.042126503159724216
$ws.Range("B24").Value = 0.041999999592715831
$ws.Range("A25").Value = -0.025167677494092544
$ws.Range("B25").Value = 0.025149055225714534
$ws.Range("A26").Value = -0.01914905549512369
$ws.Range("B26").Value = 0.01912960229109828
$ws.Range("A27").Value = -0.013129602560863596
$ws.Range("B27").Value = 0.013073315647926798
$ws.Range("A28").Value = -0.0070733159186442407
$ws.Range("B28").Value = 0.0070423162640649295
$ws.Range("A29").Value = 0.0049576834415425708
$ws.Range("B29").Value = -0.0049664576074448519
$ws.Range("A30").Value = 0.024966457282178833
$ws.Range("B30").Value = -0.025131262631942874
$ws.Range("A31").Value = 0.040131262327660622
$ws.Range("B31").Value = -0.040216663793856355
$ws.Range("A32").Value = 0.061216663467074994
$ws.Range("B32").Value = -0.061383653262995885
